$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.230.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "'2.968.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.45%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'569.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").Value = "'123.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.92%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'2.965.81"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.46%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("E10").Value = "  -6.36%  "
$ws.Range("D11").Value = "'5.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").Value = "'0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.20%  "
$ws.Range("E13").Value = "  -6.21%  "
$ws.Range("D14").Value = "'32.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.26%  "
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "'3.458.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.29%  "
$ws.Range("D17").Value = "'60.228.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "'2.968.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.30%  "
$ws.Range("D19").Value = "'6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.35%  "
$ws.Range("D20").Value = "'425.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.82%  "
$ws.Range("E21").Value = "  -6.64%  "
$ws.Range("D22").Value = "'0.658"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.70%  "
$ws.Range("E23").Value = "  -7.86%  "
$ws.Range("D24").Value = "'12.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.37%  "
$ws.Range("D25").Value = "'78.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.47%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'2.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.26%  "
$ws.Range("E29").Value = "  -7.89%  "
$ws.Range("E30").Value = "  -8.51%  "
$ws.Range("D31").Value = "'25.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.52%  "
$ws.Range("D32").Value = "'6.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.23%  "
$ws.Range("D33").Value = "'0.0920"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.11%  "
$ws.Range("D34").Value = "'2.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.71%  "
$ws.Range("E35").Value = "  -8.94%  "
$ws.Range("E36").Value = "  -4.80%  "
$ws.Range("D37").Value = "'49.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.04%  "
$ws.Range("D38").Value = "'0.0₃0647"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.10%  "
$ws.Range("D39").Value = "'7.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("E40").Value = "  -8.48%  "
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("D42").Value = "'376.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("D43").Value = "'2.628.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.93%  "
$ws.Range("E44").Value = "  -9.28%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  -6.90%  "
$ws.Range("D47").Value = "'119.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("E48").Value = "  -7.66%  "
$ws.Range("E49").Value = "  -4.90%  "
$ws.Range("D50").Value = "'23.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.94%  "
$ws.Range("D51").Value = "'31.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.52%  "
